# Automatische test-sync: 2025-06-19 13:00:10
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append the new incoming mail as row 11 on the Logs sheet ---
$logs.Range("A11").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D11").Value = "Bestelling"
$logs.Range("F11").Value = "2025-06-19 12:58:10"
$logs.Range("G11").Value = "Nee"

# --- Extend the conditional formatting ranges to cover the new row ---
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))

# --- Update the Dashboard summary table with the new "Bestelling" category ---
$dashboard.Range("A6").Value = "Bestelling"
$dashboard.Range("B6").Value = 1

# --- Extend the chart's category/value series to include the new row ---
$chart = $dashboard.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$6,'Dashboard'!`$B`$2:`$B`$6,1)"
